# Update countries & provincias Spain
# - Update the "last updated" timestamp string.
# - Armenia overtakes Costa Rica in total cases -> rows 59/60 swap countries.
# - Hungria overtakes Mauritania in total cases -> rows 104/105 swap countries.
# - Refresh the numeric stats (Casos totales, Nuevos casos, Casos activos,
#   Recuperados, Casos criticos, Muertes hoy, Muertes) for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp footer
$ws.Range("A1").Value = "Datos actualizados a 4 de Septiembre de 2020 a las 09:59"

# Rusia (row 7) - refreshed counts
$ws.Range("B7").Value = 1015105
$ws.Range("C7").Value = 5110
$ws.Range("D7").Value = 832747
$ws.Range("E7").Value = 164709
$ws.Range("G7").Value = 121
$ws.Range("H7").Value = 17649

# Armenia moves up to row 59 (was Costa Rica), Costa Rica drops to row 60
$ws.Range("A59").Value = "Armenia"
$ws.Range("B59").Value = 44461
$ws.Range("C59").Value = 190
$ws.Range("D59").Value = 39257
$ws.Range("E59").Value = 4313
$ws.Range("G59").Value = 4
$ws.Range("H59").Value = 891

$ws.Range("A60").Value = "Costa Rica"
$ws.Range("B60").Value = 44458
$ws.Range("C60").Value = 0
$ws.Range("D60").Value = 17855
$ws.Range("E60").Value = 26143
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 460

# El Salvador (row 74) - refreshed counts
$ws.Range("B74").Value = 26099
$ws.Range("C74").Value = 99
$ws.Range("D74").Value = 15347
$ws.Range("E74").Value = 10008

# Hungria moves up to row 104 (was Mauritania), Mauritania drops to row 105
$ws.Range("A104").Value = "Hungria"
$ws.Range("B104").Value = 7382
$ws.Range("C104").Value = 459
$ws.Range("D104").Value = 3944
$ws.Range("E104").Value = 2817
$ws.Range("G104").Value = 1
$ws.Range("H104").Value = 621

$ws.Range("A105").Value = "Mauritania"
$ws.Range("B105").Value = 7106
$ws.Range("C105").Value = 0
$ws.Range("D105").Value = 6588
$ws.Range("E105").Value = 358
$ws.Range("F105").Value = 0
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = 160

# Estonia (row 137) - refreshed counts
$ws.Range("B137").Value = 2456
$ws.Range("C137").Value = 15
$ws.Range("D137").Value = 2157
$ws.Range("E137").Value = 235
